# Lab1Rubric_CS295N.xlsx - "Updated the notes and lab assignment"
#
# Both worksheets ("Rubric" and "Grade") gain a new rubric line item,
# "Correct .NET Version" (worth 4 points), inserted right after the
# "MVC site" section header (i.e. becomes the new row 8, pushing the
# "GitHub repository" section and the Total row down by one row).
# At the same time every other per-item score in the "GitHub
# repository" section drops from 5 to 4 points, so the grand Total
# (SUM of the possible-points column) stays 50: 20 + 5 + 5 + 4 (new)
# + 4*4 (four items, was 4*5) = 50. Section headers ("MVC site" /
# "GitHub repository") and the "Total" row become bold, matching the
# existing bold "Lab 1" title styling.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Rubric")
$ws2 = $wb.Worksheets.Item("Grade")

# ---------------------------------------------------------------
# 1. Insert the new rubric row (row 8) on both sheets. Inserting a
#    whole row shifts everything below it down by one and Excel
#    automatically re-points the Total formulas (SUM(B4:B13) becomes
#    SUM(B4:B14), etc.) and the dimension/used-range grows to C16.
# ---------------------------------------------------------------
$ws1.Rows.Item(8).Insert()
$ws2.Rows.Item(8).Insert()

# New line item text + possible/actual scores.
$ws1.Range("A8").Value = "Correct .NET Version"
$ws1.Range("B8").Value = 4

$ws2.Range("A8").Value = "Correct .NET Version"
$ws2.Range("B8").Value = 4
$ws2.Range("C8").Value = 4

# ---------------------------------------------------------------
# 2. The rest of the "GitHub repository" section drops from 5 to 4
#    points per item (now rows 11-14 after the insert above).
# ---------------------------------------------------------------
$ws1.Range("B11:B14").Value = 4
$ws2.Range("B11:B14").Value = 4
$ws2.Range("C11:C14").Value = 4

# ---------------------------------------------------------------
# 3. Bold the section headers ("MVC site" row, "GitHub repository"
#    row) and the Total row/values, matching the bold "Lab 1" title.
#    A couple of those cells were italic before, so clear Italic
#    explicitly too (otherwise Excel keeps it and you get a
#    bold+italic font instead of plain bold).
# ---------------------------------------------------------------
foreach ($addr in @("A4", "A10", "A16", "B16")) {
    $ws1.Range($addr).Font.Bold = $true
    $ws1.Range($addr).Font.Italic = $false
}

foreach ($addr in @("A4", "A10", "A16", "B16", "C16")) {
    $ws2.Range($addr).Font.Bold = $true
    $ws2.Range($addr).Font.Italic = $false
}

# ---------------------------------------------------------------
# 4. Column widths tweaked slightly on both sheets; "Grade" also
#    gains an explicit custom width for column A (it only inherited
#    the default width before).
# ---------------------------------------------------------------
$ws1.Columns.Item(1).ColumnWidth = 24.3
$ws2.Columns.Item(1).ColumnWidth = 22.0

# ---------------------------------------------------------------
# 5. "Grade" sheet is now printed in portrait orientation.
# ---------------------------------------------------------------
$ws2.PageSetup.Orientation = 1

# ---------------------------------------------------------------
# 6. Selection / active-sheet bookkeeping: the "Rubric" tab becomes
#    the active/selected tab (it wasn't before), with the cursor
#    resting just past the bottom of its table; "Grade" keeps the
#    selection over its score table instead of being the active tab.
# ---------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A4:B16").Select() | Out-Null

$ws1.Activate()
$ws1.Range("A21").Select() | Out-Null
